# Update the "Controllers" sheet of the pyPlotList "Time series.xlsx" workbook.
# Row 3 (the TS1 controller entry) gets its yScaler / Properties / index columns
# extended so the simulation now samples two channels instead of one, making the
# start/stop time options work for both entries.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Controllers")

# D3 = yScaler, E3 = Properties, F3 = index (set in this order so the shared
# string table lands Properties, yScaler, index - matching the source order).
$ws1.Range("E3").Value = "['v.Powers','v.Powers']"
$ws1.Range("D3").Value = "[-1, -1]"
$ws1.Range("F3").Value = "['SumEven','SumOdd']"

# Leave the selection where the author ended up after editing row 3.
$ws1.Range("F11").Select()
